$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two missing Código/Transición rows (T07 / T08) that were
# previously blank on rows 8 and 9, columns P and Q.
$ws.Range("P8").Value = "T07"
$ws.Range("Q8").Value = "Plan de Iteración"
$ws.Range("P9").Value = "T08"
$ws.Range("Q9").Value = "Manual de Instalación Testify"

# Update the view so the newly filled-in columns are visible and the
# active selection reflects the last edited cell.
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("R10").Select()
